$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.393692374229431
$ws.Range("B1").Value = 2.589281797409058
$ws.Range("C1").Value = 6.662294864654541
$ws.Range("D1").Value = 2.417054653167725
$ws.Range("E1").Value = 1.202142357826233
